{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n// separates them from the \"Paulo, 1994.\" bibliography entry above them.\n\nconst jupiterResults = context.document.body.search(\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  { matchCase: true }\n);\nconst copyrightResults = context.document.body.search(\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n  { matchCase: true }\n);\njupiterResults.load(\"items\");\ncopyrightResults.load(\"items\");\nawait context.sync();\n\nif (jupiterResults.items.length === 0 || copyrightResults.items.length === 0) {\n  throw new Error(\"Expected paragraphs not found in document body.\");\n}\n\nconst jupiterParagraph = jupiterResults.items[0].paragraphs.getFirst();\nconst copyrightParagraph = copyrightResults.items[0].paragraphs.getFirst();\nconst blankParagraph = jupiterParagraph.getPrevious();\n\n// Delete blank-line + \"Ver no Jupiter...\" + \"\u00a9 2020...\" paragraphs.\ncopyrightParagraph.delete();\njupiterParagraph.delete();\nblankParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n# separates them from the \"Paulo, 1994.\" bibliography entry above them.\n\n$d = $word.ActiveDocument\n\n# Locate and delete the copyright/footer paragraph.\n$copyrightRange = $d.Content\n$copyrightFound = $copyrightRange.Find.Execute(\"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\")\nif ($copyrightFound) {\n    $copyrightParagraph = $copyrightRange.Paragraphs(1)\n    $copyrightParagraph.Range.Delete()\n}\n\n# Locate the \"Ver no Jupiter...\" paragraph, remember the blank paragraph\n# right before it, then delete both.\n$jupiterRange = $d.Content\n$jupiterFound = $jupiterRange.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif ($jupiterFound) {\n    $jupiterParagraph = $jupiterRange.Paragraphs(1)\n    $blankParagraph = $jupiterParagraph.Previous()\n    $jupiterParagraph.Range.Delete()\n    $blankParagraph.Range.Delete()\n}\n"}
